$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.091.76'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.304.55'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''300.66'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '''97.72'
$ws.Range('E6').Value = '  -2.84%  '
$ws.Range('D7').Value = '''0.520'
$ws.Range('E7').Value = '  +3.67%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.515'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '''35.73'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '''18.00'
$ws.Range('E13').Value = '  -3.10%  '
$ws.Range('D14').Value = '''6.87'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').Value = '2.662.38'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '2.349.98'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '''0.786'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '42.977.08'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '''13.26'
$ws.Range('E19').Value = '  +7.45%  '
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('D21').Value = '''6.11'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = '''68.33'
$ws.Range('E22').Value = '  +0.42%  '
$ws.Range('D23').Value = '''238.37'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').Value = '''2.20'
$ws.Range('E24').Value = '  -2.93%  '
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '''2.42'
$ws.Range('E26').Value = '  -1.75%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''24.74'
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''167.67'
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''9.15'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.04'
$ws.Range('E30').Value = '  -12.46%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '''32.70'
$ws.Range('E31').Value = '  -5.74%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '''0.998'
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''5.15'
$ws.Range('E33').Value = '  +2.37%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = '''18.12'
$ws.Range('E34').Value = '  +1.57%  '
$ws.Range('D35').Value = '''4.79'
$ws.Range('E35').Value = '  +1.74%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''2.41'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.0688'
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '''0.102'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''1.79'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '''0.111'
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '''2.75'
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.009.74'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0288'
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = '''2.16'
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''10.17'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''17.19'
$ws.Range('E46').Value = '  -2.89%  '
$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').Value = '''54.34'
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.524.99'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '''1.53'
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = '''2.80'
$ws.Range('E51').Value = '  +10.78%  '
